# remove Gamelogic project, modify SLG building config
#
# Row 7 (A7/J7) used to describe "Desc" (描述). It is repurposed to describe
# "Icon" (图标), and two new rows are appended:
#   Row 8: ShowName / 名字
#   Row 9: Desc / 描述 (the original row-7 content, now relocated)
#
# Columns: A=Id, B=Type, C=Public, D=Private, E=Save, F=View, G=Index,
#          H=SaveInterval, I=RelationValue, J=Desc

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: repurpose the old "Desc" row into the "Icon" row -----------
$ws.Cells.Item(7, 1).Value = "Icon"
$ws.Cells.Item(7, 10).Value = "图标"

# --- Row 8: new "ShowName" row ------------------------------------------
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "ShowName"

$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 2).Value = "string"

$ws.Cells.Item(8, 3).Value = $False
$ws.Cells.Item(8, 4).Value = $False
$ws.Cells.Item(8, 5).Value = $False
$ws.Cells.Item(8, 6).Value = $False
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0

$ws.Cells.Item(8, 9).NumberFormat = "@"
$ws.Cells.Item(8, 9).Value = "Friend"

$ws.Cells.Item(8, 10).NumberFormat = "@"
$ws.Cells.Item(8, 10).Value = "名字"

# --- Row 9: relocated "Desc" row (what used to live in row 7) ----------
$ws.Cells.Item(9, 1).NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "Desc"

$ws.Cells.Item(9, 2).NumberFormat = "@"
$ws.Cells.Item(9, 2).Value = "string"

$ws.Cells.Item(9, 3).Value = $False
$ws.Cells.Item(9, 4).Value = $False
$ws.Cells.Item(9, 5).Value = $False
$ws.Cells.Item(9, 6).Value = $False
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0

$ws.Cells.Item(9, 9).NumberFormat = "@"
$ws.Cells.Item(9, 9).Value = "Friend"

$ws.Cells.Item(9, 10).NumberFormat = "@"
$ws.Cells.Item(9, 10).Value = "描述"

# --- Data validation list for column F now starts after the new rows ---
$oldRange = $ws.Range("F8:F1048576")
$oldRange.Validation.Delete()
$newRange = $ws.Range("F10:F1048576")
$newRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Reflect the saved selection state in the sheet view ---------------
$ws.Range("C13").Select()
